$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark D5 as "worked" like C5 (same style), which bumps the totals formulas
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D5").Value = 1

# Move the active selection to B10 (single cell)
$ws.Range("B10").Select()
